# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (positioned right before "总计")
#   with the quarter's fund-holding detail rows.
# - Updates the "总计" (totals) sheet: adds a new top row for 2022-Q1
#   and shifts the existing 2021-Q4 / 2021-Q3 summary rows down,
#   renumbering the index column.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted immediately before "总计"
#    (NOTE: grab the "总计" handle, use it once as the Add() anchor, then
#    re-resolve "总计" by name afterwards -- the handle captured before
#    Add() tracks the *position*, which the insert shifts by one.)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# Borrow the header-row / index-column formatting already used by the
# other quarterly sheets (e.g. 2021-Q4) instead of re-building styles.
$q4Sheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:A4").Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$fundRows = @(
  @(0,  "009794", "太平智选一年定期开放股票",       "5.20", "89.96", "4.58", "0.2382", 6),
  @(1,  "013422", "太平智行三个月定期开放混合",      "6.80", "81.87", "3.03", "0.2060", 8),
  @(2,  "002293", "南方益和灵活配置混合",            "1.40", "83.59", "4.34", "0.0608", 3),
  @(3,  "011471", "鹏华致远成长混合A",               "2.19", "61.03", "2.62", "0.0574", 8),
  @(4,  "004194", "招商中证1000指数增强A",           "1.76", "94.40", "1.08", "0.0190", 6),
  @(5,  "006522", "财通新兴蓝筹混合A",               "0.29", "90.33", "3.65", "0.0106", 9),
  @(6,  "004195", "招商中证1000指数增强C",           "0.68", "94.40", "1.08", "0.0073", 6),
  @(7,  "970073", "东证融汇成长优选混合A",           "0.68", "82.02", "0.87", "0.0059", 4),
  @(8,  "970074", "东证融汇成长优选混合C",           "0.27", "82.02", "0.87", "0.0023", 4),
  @(9,  "011472", "鹏华致远成长混合C",               "0.08", "61.03", "2.62", "0.0021", 8),
  @(10, "006523", "财通新兴蓝筹混合C",               "0.03", "90.33", "3.65", "0.0011", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $ws.Range("A$r").Value = $row[0]

    # Force text on the numeric-looking fields (fund code / AUM / weight
    # figures) so leading zeros and exact decimal text survive, same as
    # they're stored on every other quarterly sheet.
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $row[1]

    $ws.Range("C$r").Value = $row[2]

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $row[3]

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]

    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = $row[5]

    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("G$r").Value = $row[6]

    $ws.Range("H$r").Value = $row[7]

    $r++
}

# Now that every data cell has its real value, strip the temporary "@"
# text format back out (paste-format-only from the already-unstyled
# name column) so the data rows end up on the default style, matching
# the other quarterly sheets.
$ws.Range("C2").Copy()
$ws.Range("B2:G12").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q1 summary row on top, push the rest down
#    (re-resolved by name now that the sheet order has settled)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 0.61

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Re-apply the index-column style (s=2) to the newly inserted A2 cell
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
